$wb = $excel.ActiveWorkbook

# --- Sheet "Reguły" (8th sheet) : reorder object lists inside a few rule descriptions ---
$rules = $wb.Worksheets.Item(8)

$rules.Range("B2").Value = "(attempts >=  3.0) & (pregnancy <=  0.0) => (class <= 1) ['a1', 'a3', 'a7']"
$rules.Range("B6").Value = "(age >=  42.0) => (class <= 1) ['a3', 'a14']"
$rules.Range("B7").Value = "(age <=  31.0) & (attempts <=  1.0) & (endometrium <=  1.0) => (class >= 2) ['a11', 'a24', 'a12', 'a9']"
$rules.Range("B8").Value = "(frozen_embryos >=  8.0) & (sperm <=  1.0) => (class >= 2) ['a6', 'a16']"

# --- Sheet "Walidacja krzyżowa" (10th sheet) : add experiment results (reordered metrics) ---
$cv = $wb.Worksheets.Item(10)

$cv.Range("A1").Value = "correct"
$cv.Range("B1").Value = 0.6923076923076923

$cv.Range("A2").Value = "not_classified"
$cv.Range("B2").Value = 0.48

$cv.Range("A3").Value = "f1_score"
$cv.Range("B3").Value = 0.4466666666666665

$cv.Range("A4").Value = "accuracy"
$cv.Range("B4").Value = 0.36
